$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells in existing rows 2-11 (columns E, F, H, I, K, L)
$ws.Range("E2").Value = -17.14
$ws.Range("F2").Value = 73.41
$ws.Range("H2").Value = 3.7
$ws.Range("I2").Value = 67.77
$ws.Range("K2").Value = -0.94
$ws.Range("L2").Value = 73.77

$ws.Range("E3").Value = 12.31
$ws.Range("F3").Value = 57.35
$ws.Range("H3").Value = 11.64
$ws.Range("I3").Value = 60.94
$ws.Range("K3").Value = 27.34
$ws.Range("L3").Value = 57.29

$ws.Range("E4").Value = 11.85
$ws.Range("F4").Value = 53.08
$ws.Range("H4").Value = 22.11
$ws.Range("I4").Value = 53.22
$ws.Range("K4").Value = 16.25
$ws.Range("L4").Value = 58.57

$ws.Range("E5").Value = -6.4
$ws.Range("F5").Value = 64.58
$ws.Range("H5").Value = 8.76
$ws.Range("I5").Value = 64.2
$ws.Range("K5").Value = 5.3
$ws.Range("L5").Value = 66.45

$ws.Range("E6").Value = 12.84
$ws.Range("F6").Value = 45.2
$ws.Range("H6").Value = 16.11
$ws.Range("I6").Value = 51.24
$ws.Range("K6").Value = 49.07
$ws.Range("L6").Value = 31.29

$ws.Range("E7").Value = 13.25
$ws.Range("F7").Value = 53.45
$ws.Range("H7").Value = 32.76
$ws.Range("I7").Value = 47.84
$ws.Range("K7").Value = 28.45
$ws.Range("L7").Value = 51.29

$ws.Range("E8").Value = 35.63
$ws.Range("F8").Value = 39.03
$ws.Range("H8").Value = 20.6
$ws.Range("I8").Value = 54.27
$ws.Range("K8").Value = 31.95
$ws.Range("L8").Value = 50.77

$ws.Range("E9").Value = 6.86
$ws.Range("F9").Value = 54.32
$ws.Range("H9").Value = 23.53
$ws.Range("I9").Value = 51.7
$ws.Range("K9").Value = 20.07
$ws.Range("L9").Value = 56.23

$ws.Range("E10").Value = 52.19
$ws.Range("F10").Value = 21.65
$ws.Range("H10").Value = 74.67
$ws.Range("I10").Value = 14.85
$ws.Range("K10").Value = 77.6
$ws.Range("L10").Value = 13.37

$ws.Range("E11").Value = -3.12
$ws.Range("F11").Value = 53.33
$ws.Range("H11").Value = 17.91
$ws.Range("I11").Value = 52.99
$ws.Range("K11").Value = 10.83
$ws.Range("L11").Value = 58.13

# Add new rows 12-21
$ws.Range("A12").Value = 4.14
$ws.Range("B12").Value = "Planar"
$ws.Range("C12").Value = 2.5
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1.31
$ws.Range("F12").Value = 402.84
$ws.Range("G12").Value = 408.21
$ws.Range("H12").Value = -0.03
$ws.Range("I12").Value = 403.28
$ws.Range("J12").Value = 403.14
$ws.Range("K12").Value = 71.57
$ws.Range("L12").Value = 120.84
$ws.Range("M12").Value = 424.88

$ws.Range("A13").Value = 2.5
$ws.Range("B13").Value = "Planar"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 2.5
$ws.Range("E13").Value = 63.25
$ws.Range("F13").Value = 216.48
$ws.Range("G13").Value = 589.2
$ws.Range("H13").Value = 0.18
$ws.Range("I13").Value = 577.81
$ws.Range("J13").Value = 578.8
$ws.Range("K13").Value = 65.25
$ws.Range("L13").Value = 215.85
$ws.Range("M13").Value = 621.2

$ws.Range("A14").Value = 4.8
$ws.Range("B14").Value = "Planar"
$ws.Range("C14").Value = 2.5
$ws.Range("D14").Value = 2.5
$ws.Range("E14").Value = 1.44
$ws.Range("F14").Value = 336.94
$ws.Range("G14").Value = 341.67
$ws.Range("H14").Value = 0.18
$ws.Range("I14").Value = 337.55
$ws.Range("J14").Value = 337.92
$ws.Range("K14").Value = 65.77
$ws.Range("L14").Value = 123.5
$ws.Range("M14").Value = 360.62

$ws.Range("A15").Value = 16.08
$ws.Range("B15").Value = "Planar"
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = -1.37
$ws.Range("F15").Value = 205.21
$ws.Range("G15").Value = 202.43
$ws.Range("H15").Value = 0.04
$ws.Range("I15").Value = 204.99
$ws.Range("J15").Value = 205.04
$ws.Range("K15").Value = 62.76
$ws.Range("L15").Value = 77.57
$ws.Range("M15").Value = 208.27

$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Planar"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 2.19
$ws.Range("F16").Value = 175.59
$ws.Range("G16").Value = 179.5
$ws.Range("H16").Value = 0.93
$ws.Range("I16").Value = 178.43
$ws.Range("J16").Value = 180.1
$ws.Range("K16").Value = 6.6
$ws.Range("L16").Value = 178.23
$ws.Range("M16").Value = 190.8

$ws.Range("A17").Value = 20.04
$ws.Range("B17").Value = "Planar"
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 59.39
$ws.Range("F17").Value = 83.85
$ws.Range("G17").Value = 206.49
$ws.Range("H17").Value = 0.3
$ws.Range("I17").Value = 209.08
$ws.Range("J17").Value = 209.73
$ws.Range("K17").Value = 4.07
$ws.Range("L17").Value = 202.82
$ws.Range("M17").Value = 211.43

$ws.Range("A18").Value = 41.2
$ws.Range("B18").Value = "Gradual"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 2.5
$ws.Range("E18").Value = 42.38
$ws.Range("F18").Value = 83.06
$ws.Range("G18").Value = 144.15
$ws.Range("H18").Value = 1.87
$ws.Range("I18").Value = 148.1
$ws.Range("J18").Value = 150.9
$ws.Range("K18").Value = 4.72
$ws.Range("L18").Value = 145.07
$ws.Range("M18").Value = 152.23

$ws.Range("A19").Value = 43.99
$ws.Range("B19").Value = "Gradual"
$ws.Range("C19").Value = 2.5
$ws.Range("D19").Value = 2.5
$ws.Range("E19").Value = 46.55
$ws.Range("F19").Value = 70.75
$ws.Range("G19").Value = 132.35
$ws.Range("H19").Value = 6.09
$ws.Range("I19").Value = 132.04
$ws.Range("J19").Value = 140.6
$ws.Range("K19").Value = 4.5
$ws.Range("L19").Value = 135.63
$ws.Range("M19").Value = 142.01

$ws.Range("A20").Value = 71.88
$ws.Range("B20").Value = "Gradual"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 58.73
$ws.Range("F20").Value = 47.75
$ws.Range("G20").Value = 115.71
$ws.Range("H20").Value = 61.23
$ws.Range("I20").Value = 48.21
$ws.Range("J20").Value = 124.36
$ws.Range("K20").Value = 4.84
$ws.Range("L20").Value = 119.17
$ws.Range("M20").Value = 125.24

$ws.Range("A21").Value = 83.96
$ws.Range("B21").Value = "Gradual"
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 58.51
$ws.Range("F21").Value = 49.84
$ws.Range("G21").Value = 120.1
$ws.Range("H21").Value = 60.17
$ws.Range("I21").Value = 51.62
$ws.Range("J21").Value = 129.59
$ws.Range("K21").Value = 4.38
$ws.Range("L21").Value = 123.85
$ws.Range("M21").Value = 129.51
